# Update F2, F3, F4 values in both "展览" and "全部类型" worksheets
# F2: 4900 -> 4920
# F3: 146  -> 148
# F4: 856  -> 859

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 4920
    $ws.Range("F3").Value = 148
    $ws.Range("F4").Value = 859
}
